$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$row = 47

$ws.Range("A$row").Value = "Staan er nog EcoPro-700 op voorraad?"
$ws.Range("B$row").Value = "mailmind.test@zohomail.eu"
$ws.Range("C$row").Value = "Testmail #6: Staan er nog EcoPro-700 op voorraad?"
$ws.Range("D$row").Value = "Overig"
$ws.Range("E$row").Value = "Beste klant,`nDank u voor uw e-mail. Op dit moment hebben we nog EcoPro-700 op voorraad. Als u geïnteresseerd bent in het plaatsen van een bestelling, kunt u dit via onze website doen of contact opnemen met onze verkoopafdeling. Mocht u verdere vragen hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$ws.Range("F$row").Value = "2025-08-05 19:33:16"
$ws.Range("G$row").Value = "Ja"
$ws.Range("H$row").Value = "Nee"
$ws.Range("I$row").Value = "Ja"
$ws.Range("J$row").Value = "Nee"

# Avoid the engine stamping an explicit custom row height on the
# newly-created row - re-measure it like the rest of the sheet.
$ws.Rows.Item($row).EntireRow.AutoFit()

# Extend the conditional-formatting ranges so they cover the new row,
# matching how Excel grows these when new rows are appended.
$ws.Range("D2:D46").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D47"))
$ws.Range("G2:G46").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G47"))
$ws.Range("H2:H46").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H47"))
$ws.Range("I2:I46").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I47"))
$ws.Range("J2:J46").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J47"))

# Update the Dashboard "Overig" category count (row 5) from 4 to 5.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B5").Value = 5
